$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.668.95'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '1.643.75'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("E6").Value = '  +1.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +1.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0841'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").Value = '1.872.62'
$ws.Range("E13").Value = '  +2.19%  '
$ws.Range("D14").Value = '1.630.43'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("E15").Value = '  +1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.64%  '
$ws.Range("D17").Value = '26.747.89'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").Value = '0.0₃0745'
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("E24").Value = '  +15.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("E28").Value = '  +4.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0517'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.65%  '
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("E33").Value = '  +3.72%  '
$ws.Range("D34").Value = '1.275.82'
$ws.Range("E34").Value = '  +4.93%  '
$ws.Range("E35").Value = '  +2.57%  '
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0181'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.534'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.833'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.34%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.816'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.74%  '
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.13%  '
$ws.Range("D44").Value = '1.783.03'
$ws.Range("E44").Value = '  +1.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.61'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0516'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("E49").Value = '  +3.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0968'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("E51").Value = '  -0.45%  '
